$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the containing path prefix from the fastq filenames in column F
$ws.Range("F2").Value = "s_6_sequence_30C.fastq.gz"
$ws.Range("F3").Value = "s_6_sequence_37C_CO2.fastq.gz"
$ws.Range("F4").Value = "s_6_sequence_CHX_30C.fastq.gz"
$ws.Range("F5").Value = "s_6_sequence_CHX_37C_CO2.fastq.gz"

# Update the active cell selection to F5
$ws.Range("F5").Select()
